$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6

# Row 3 updates
$ws.Range("G3").Value = 8.5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 1.36
$ws.Range("K3").Value = 2.38
$ws.Range("S3").Value = 1.36
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
